# Applies the "bia_faim_models" schema restructuring:
#  - Drop the anyOf/oneOf wrapper sheets (OrganisationInfoCollection,
#    OrganisationURL, FundingStatement, GrantReferenceCollection)
#  - Fold their single field back into the sheet that used them
#  - Add a new FileLevelMetadata class/sheet
#  - Add a second example slot pairing (Study.grants, Annotations.file_metadata)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update "Study": funding -> funding_statement, add grants column
# ---------------------------------------------------------------------------
$study = $wb.Worksheets.Item("Study")
$study.Cells.Item(1, 7).Value = "funding_statement"
$study.Cells.Item(1, 12).Value = "grants"

# ---------------------------------------------------------------------------
# 2. Update "Author": affiliation/role -> role/organisation
# ---------------------------------------------------------------------------
$author = $wb.Worksheets.Item("Author")
$author.Cells.Item(1, 5).Value = "role"
$author.Cells.Item(1, 6).Value = "organisation"

# ---------------------------------------------------------------------------
# 3. Fold OrganisationURL.ror_id into OrganisationInfo, then drop the
#    now-empty collection/url wrapper sheets
# ---------------------------------------------------------------------------
$orgInfo = $wb.Worksheets.Item("OrganisationInfo")
$orgInfo.Cells.Item(1, 3).Value = "ror_id"

$wb.Worksheets.Item("OrganisationInfoCollection").Delete()
$wb.Worksheets.Item("OrganisationURL").Delete()

# ---------------------------------------------------------------------------
# 4. Fold FundingStatement back into Study.funding_statement (already set
#    above) and drop the wrapper sheet
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("FundingStatement").Delete()

# ---------------------------------------------------------------------------
# 5. Drop GrantReferenceCollection (its "grants" field now lives on Study)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("GrantReferenceCollection").Delete()

# ---------------------------------------------------------------------------
# 6. Insert the new "FileLevelMetadata" class right after GrantReference
# ---------------------------------------------------------------------------
$afterGrantReference = $wb.Worksheets.Item("GrantReference")
$fileMeta = $wb.Worksheets.Add($null, $afterGrantReference)
$fileMeta.Name = "FileLevelMetadata"

$fileMeta.Cells.Item(1, 1).Value = "annotation_id"
$fileMeta.Cells.Item(1, 2).Value = "annotation_type"
$fileMeta.Cells.Item(1, 3).Value = "source_image_id"
$fileMeta.Cells.Item(1, 4).Value = "transformations"
$fileMeta.Cells.Item(1, 5).Value = "spatial_information"

$fileMetaDv = $fileMeta.Range("B2:B1048576").Validation
$fileMetaDv.Add(3, 1, 1, """class_labels,bounding_boxes,counts,derived_annotations,geometrical_annotations,graphs,point_annotations,segmentation_mask,tracks,weak_annotations,other""")

# ---------------------------------------------------------------------------
# 7. Rewrite "Annotations": drop annotation_type/source_image/
#    spatial_information/transformatons columns, add file_metadata, and
#    remove the annotation_type list validation (it now lives on
#    FileLevelMetadata.annotation_type)
# ---------------------------------------------------------------------------
$annotations = $wb.Worksheets.Item("Annotations")
$annotationsDv = $annotations.Range("B2:B1048576").Validation
$annotationsDv.Delete()

$annotations.Cells.Item(1, 1).Value = "annotation_overview"
$annotations.Cells.Item(1, 2).Value = "annotation_method"
$annotations.Cells.Item(1, 3).Value = "annotation_criteria"
$annotations.Cells.Item(1, 4).Value = "annotation_coverage"
$annotations.Cells.Item(1, 5).Value = "annotation_confidence_level"
$annotations.Cells.Item(1, 6).Value = "authors"
$annotations.Cells.Item(1, 7).Value = "file_metadata"
$annotations.Cells.Item(1, 8).Value = ""
$annotations.Cells.Item(1, 9).Value = ""
$annotations.Cells.Item(1, 10).Value = ""

# ---------------------------------------------------------------------------
# 8. "Version" and "GrantReference" content is unchanged - only their
#    sheetId/position shifts, which naturally follows from the deletes above.
# ---------------------------------------------------------------------------
